$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "2009年" data row), shifting subsequent rows up
$ws.Rows.Item(2).Delete()
